# Applies the edits described by the diff: updates the date line and the
# division problems throughout the tables.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-02-27 Thursday" "2025-02-28 Friday"

Replace-Text "34÷5=" "68÷6="
Replace-Text "85÷2=" "47÷2="
Replace-Text "60÷6=" "83÷9="
Replace-Text "42÷7=" "97÷7="
Replace-Text "17÷9=" "49÷9="
Replace-Text "39÷2=" "53÷5="
Replace-Text "99÷7=" "46÷9="
Replace-Text "71÷7=" "69÷4="
Replace-Text "46÷7=" "58÷4="
Replace-Text "51÷2=" "93÷5="
Replace-Text "59÷2=" "70÷9="
Replace-Text "96÷7=" "46÷9="
Replace-Text "75÷4=" "40÷2="
Replace-Text "29÷3=" "23÷4="
Replace-Text "80÷8=" "83÷6="
Replace-Text "57÷9=" "99÷5="
Replace-Text "20÷8=" "76÷7="
Replace-Text "78÷7=" "57÷6="
Replace-Text "73÷3=" "14÷2="
Replace-Text "97÷5=" "56÷8="
Replace-Text "27÷3=" "11÷2="
Replace-Text "29÷6=" "27÷9="
Replace-Text "80÷6=" "45÷6="
Replace-Text "86÷3=" "34÷8="
Replace-Text "30÷4=" "65÷2="
